# daily auto push: 2026-01-11 18:42 UTC
#
# The data table in column A jumps straight from 2026/01/11 to 2026/12/29
# (row 604 in the "before" state). This push fills in the two missing
# 2026/01/11 / 2026/01/12 rows (continuing the existing time-of-day
# sequence) by inserting two new rows right before the old row 604 and
# pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the position right before the "2026/12/29" block
# (old row 604), shifting everything from there on down by two rows.
$ws.Rows.Item(604).Insert()
$ws.Rows.Item(604).Insert()

# --- New row 604: 2026/01/11 日 23 19 -------------------------------------
# Force column A to be stored as literal text (not auto-parsed into a date
# serial number) by temporarily marking the cell as text, then clear the
# formatting residue afterwards so the cell ends up as a plain text cell,
# matching every other date cell in the column.
$ws.Range("A604").NumberFormat = "@"
$ws.Range("A604").Value = "2026/01/11"
$ws.Range("A604").ClearFormats()
$ws.Range("B604").Value = "日"
$ws.Range("C604").Value = 23
$ws.Range("D604").Value = 19

# --- New row 605: 2026/01/12 月 1 15 --------------------------------------
$ws.Range("A605").NumberFormat = "@"
$ws.Range("A605").Value = "2026/01/12"
$ws.Range("A605").ClearFormats()
$ws.Range("B605").Value = "月"
$ws.Range("C605").Value = 1
$ws.Range("D605").Value = 15
